# Add a new column I ("6/14" acumulado) to the daily-cumulative deaths sheet,
# mirroring the existing H column (formula total, header date, daily values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: running total formula (plain style, like H1) ---
$ws.Range("I1").Formula = "=SUM(I3:I93)"

# --- Row 2: new date header, copying H2's format (border+fill+date numfmt) ---
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").Value = 43996
$excel.CutCopyMode = 0

# --- Rows 3-83: daily values for the new column ---
$values = @{
    3 = 1;   4 = 1;   5 = 1;   6 = 1;   7 = 1;   8 = 1;   9 = 3;   10 = 4;
    11 = 1;  12 = 4;  13 = 6;  14 = 8;  15 = 1;  16 = 5;  17 = 4;  18 = 10;
    19 = 9;  20 = 6;  21 = 7;  22 = 3;  23 = 9;  24 = 4;  25 = 8;  26 = 6;
    27 = 5;  28 = 6;  29 = 5;  30 = 10; 31 = 8;  32 = 5;  33 = 7;  34 = 6;
    35 = 7;  36 = 9;  37 = 11; 38 = 7;  39 = 6;  40 = 8;  41 = 10; 42 = 7;
    43 = 19; 44 = 12; 45 = 11; 46 = 14; 47 = 29; 48 = 16; 49 = 22; 50 = 30;
    51 = 19; 52 = 25; 53 = 38; 54 = 35; 55 = 32; 56 = 49; 57 = 45; 58 = 60;
    59 = 60; 60 = 65; 61 = 60; 62 = 69; 63 = 85; 64 = 95; 65 = 86; 66 = 103;
    67 = 115;68 = 111;69 = 120;70 = 121;71 = 108;72 = 136;73 = 129;74 = 140;
    75 = 134;76 = 120;77 = 136;78 = 131;79 = 125;80 = 140;81 = 123;82 = 110;
    83 = 24
}

foreach ($r in $values.Keys) {
    $ws.Cells.Item($r, 9).Value = $values[$r]
}

# --- Sheet view: scroll back to top and move the active selection ---
$ws.Range("L7").Select()
